$d = $word.ActiveDocument

# 1. "Curso (semestre ideal): EQN (12)" -> "Curso (semestre ideal): EQD (10), EQN (12)"
$d.Content.Find.Execute(
    "Curso (semestre ideal): EQN (12)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Curso (semestre ideal): EQD (10), EQN (12)", 2)

# 2. Remove the "Requisitos" heading paragraph and the
#    "LOQ4044 - ... (Requisito fraco)" bullet paragraph that followed it
#    at the end of the document.
$reqPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Requisitos`r") {
        $reqPara = $p
    }
}

if ($reqPara -ne $null) {
    $r = $d.Range($reqPara.Range.Start, $d.Content.End)
    $r.Delete()
}
